$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

# Add new rows of account data
$ws.Range("A3").Value = 11700
$ws.Range("B3").Value = "Doha@2021"

$ws.Range("A4").Value = 2168
$ws.Range("B4").Value = "Qatar@2021"

$ws.Range("A5").Value = 88996
$ws.Range("B5").Value = "Qatar@2021"

$ws.Range("A6").Value = 111880
$ws.Range("B6").Value = "Qatar@2021"

$ws.Range("A7").Value = 11848
$ws.Range("B7").Value = "Qatar@2021"

$ws.Range("A8").Value = 80197
$ws.Range("B8").Value = "Qatar@2021"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Doha@2021")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Qatar@2021")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:Qatar@2021")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:Qatar@2021")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:Qatar@2021")

$ws.Range("A2:B8").Borders.LineStyle = 1
$ws.Range("A2:A8").HorizontalAlignment = -4131

$ws.Range("A5").Select()
$wb.Windows.Item(1).WindowState = -4143
